$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMID")

# Insert a new column before D, restricted to the used data range (A5:K102),
# so D:K shift right to E:L for every row without materializing the whole sheet.
$ws.Range("D5:D102").Insert(-4161)

# The newly inserted column D starts out General-formatted; restore the
# number/date formatting by copying it over from column E (the shifted former D).
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly inserted column D with the newest quarter of data
# (period ending 43373, i.e. 2018-09-30)
$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 202300
$ws.Range("D9").Value = 150300
$ws.Range("D10").Value = 52000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 17000
$ws.Range("D15").Value = 23000
$ws.Range("D17").Value = 135200
$ws.Range("D18").Value = 67100
$ws.Range("D20").Value = 24900
$ws.Range("D21").Value = 115100
$ws.Range("D22").Value = 22600
$ws.Range("D23").Value = 69400
$ws.Range("D24").Value = 31200
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 38200
$ws.Range("D27").Value = 29400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -24900
$ws.Range("D33").Value = 29400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 29400
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 22800
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 92400
$ws.Range("D44").Value = 3100
$ws.Range("D45").Value = 170800
$ws.Range("D46").Value = 289100
$ws.Range("D47").Value = 336800
$ws.Range("D48").Value = 993900
$ws.Range("D49").Value = 188300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 26600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1834600
$ws.Range("D57").Value = 48300
$ws.Range("D58").Value = 603500
$ws.Range("D59").Value = 132400
$ws.Range("D60").Value = 784300
$ws.Range("D61").Value = 501200
$ws.Range("D62").Value = 85200
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1384500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 250900
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 199200
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = 29400
$ws.Range("D83").Value = 23000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 4000
$ws.Range("D91").Value = -16800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 188200
$ws.Range("D96").Value = 52500
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -185600
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 6600

# A few historical quarters (Capital Expenditures row 91, Dividends Paid row 96)
# were also restated as part of this update
$ws.Range("H91").Value = -22000
$ws.Range("I91").Value = -44000
$ws.Range("J91").Value = -20200
$ws.Range("H96").Value = -15100
